# Revert "new changes in ops (ordercreation & orderpage & order form)"
# This reverts the data-entry / header-reorder edits back to the prior
# (template) shape: header columns F:H swap back to Lob/Process/Product Name
# order, the two sample rows get replaced with the original sample data,
# and a couple of cells that had a one-off font/border style get normalised
# back to the common data-row style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1): F/G/H swap back -------------------------------
$ws.Range("F1").Value = "Product Name"
$ws.Range("G1").Value = "Lob"
$ws.Range("H1").Value = "Process"

# ---- Normalise the one-off styles on C2:E2/C3:E3 to the common style ---
# (before: C/D/E used dedicated font+border xfs; after: same plain xf as
# the rest of the row, e.g. B2/B3.) Copy formats from B2 which already
# carries the target style.
$ws.Range("B2").Copy()
$ws.Range("C2:E2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B3").Copy()
$ws.Range("C3:E3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# ---- Row 2 data ----------------------------------------------------------
$ws.Range("A2").Value = 45436
$ws.Range("B2").Value = 1213286
$ws.Range("C2").Value = "SIPL0005"
$ws.Range("D2").Value = "SIPL0004"
$ws.Range("E2").Value = "Reltco"
$ws.Range("F2").Value = "Commercial Full Search"
$ws.Range("G2").Value = "Title"
$ws.Range("H2").Value = "Search"
$ws.Range("I2").Value = "FL"
$ws.Range("J2").Value = "Clay"
$ws.Range("K2").Value = "FLClay"
$ws.Range("L2").Value = "WIP"
$ws.Range("M2").Value = "Search(T1) "

# ---- Row 3 data ----------------------------------------------------------
$ws.Range("A3").Value = 45439
$ws.Range("B3").Value = 2193289
$ws.Range("C3").Value = "SIPL0005"
$ws.Range("D3").Value = "SIPL0004"
$ws.Range("E3").Value = "Reltco"
$ws.Range("F3").Value = "Residential Current Owner Search"
$ws.Range("G3").Value = "Title"
$ws.Range("H3").Value = "Search"
$ws.Range("I3").Value = "FL"
$ws.Range("J3").Value = "Clay"
$ws.Range("K3").Value = "FLClay"
$ws.Range("L3").Value = "WIP"
$ws.Range("M3").Value = "Search(T2)"

# ---- Column widths --------------------------------------------------------
$ws.Range("C1").ColumnWidth = 35.5
$ws.Range("F1:H1").ColumnWidth = 15.17

# ---- Selection -------------------------------------------------------------
$ws.Range("F13").Select()
